# Update the CodeSystem document-format implementation guide spreadsheet.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Date" property value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2022-05-18T17:38:26+00:00"

# --- Concepts sheet: clear the "Definition" column values (D2:D6) ---
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("D2:D6").ClearContents()
